$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 86
$ws.Range("I4").Value = 86
$ws.Range("K4").Value = 86
$ws.Range("M4").Value = 28

$ws.Range("H40").Value = 12075.15
$ws.Range("I40").Value = 1875.375
$ws.Range("K40").Value = 1875.375
$ws.Range("M40").Value = -1700.375

$ws.Range("H55").Value = 199.6
$ws.Range("I55").Value = 100
$ws.Range("J55").Value = 224.5
$ws.Range("K55").Value = 100
$ws.Range("L55").Value = 224.5
$ws.Range("M55").Value = 114
$ws.Range("N55").Value = -652.5

$ws.Range("H113").Value = 5132.6665
$ws.Range("I113").Value = 4700
$ws.Range("J113").Value = 5998
$ws.Range("K113").Value = 4700
$ws.Range("L113").Value = 5998
$ws.Range("M113").Value = -1446
$ws.Range("N113").Value = -12506

$ws.Range("H132").Value = 2369.7954
$ws.Range("I132").Value = 2331.359
$ws.Range("K132").Value = 6994.076999999999
$ws.Range("M132").Value = -4464.076999999999

$ws.Range("H137").Value = 4386.364
$ws.Range("I137").Value = 4365.625
$ws.Range("K137").Value = 13096.875
$ws.Range("M137").Value = -10546.875

$ws.Range("H138").Value = 3614.8438
$ws.Range("J138").Value = 3949.6538
$ws.Range("L138").Value = 11848.9614
$ws.Range("N138").Value = -22128.9614

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 7399.8
$ws.Range("I2").Value = 7399.8
$ws.Range("K2").Value = 7399.8
$ws.Range("M2").Value = -7286.8

$ws.Range("H74").Value = 1999.9354
$ws.Range("I74").Value = 1509.95
$ws.Range("K74").Value = 1509.95
$ws.Range("M74").Value = -635.95

$ws.Range("H77").Value = 1999.9354
$ws.Range("I77").Value = 1509.95
$ws.Range("K77").Value = 7549.75
$ws.Range("M77").Value = -3181.75

$ws.Range("H95").Value = 26463
$ws.Range("J95").Value = 26463
$ws.Range("L95").Value = 26463
$ws.Range("N95").Value = -31955

$ws.Range("H110").Value = 1793.2727
$ws.Range("I110").Value = 1797.6
$ws.Range("K110").Value = 1797.6
$ws.Range("M110").Value = 247.4000000000001

$ws.Range("H116").Value = 7399.8
$ws.Range("I116").Value = 7399.8
$ws.Range("K116").Value = 7399.8
$ws.Range("M116").Value = -5105.8

$ws.Range("H122").Value = 3230.6843
$ws.Range("J122").Value = 3536.5
$ws.Range("L122").Value = 10609.5
$ws.Range("N122").Value = -15509.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 7399.8
$ws.Range("I3").Value = 7399.8
$ws.Range("K3").Value = 7399.8
$ws.Range("M3").Value = -7285.8

$ws.Range("H99").Value = 4223.077
$ws.Range("I99").Value = 3990.9092
$ws.Range("K99").Value = 3990.9092
$ws.Range("M99").Value = -2492.9092

$ws.Range("H105").Value = 3455.3635
$ws.Range("I105").Value = 3012.375
$ws.Range("K105").Value = 3012.375
$ws.Range("M105").Value = -1265.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 78081.94
$ws.Range("I16").Value = 40908.547
$ws.Range("K16").Value = 40908.547
$ws.Range("M16").Value = -40621.547

$ws.Range("H58").Value = 5782
$ws.Range("I58").Value = 3912
$ws.Range("J58").Value = 6249.5
$ws.Range("K58").Value = 3912
$ws.Range("L58").Value = 6249.5
$ws.Range("M58").Value = -3709
$ws.Range("N58").Value = -6655.5

$ws.Range("H81").Value = 34908.668
$ws.Range("J81").Value = 34908.668
$ws.Range("L81").Value = 34908.668
$ws.Range("N81").Value = -36904.668

$ws.Range("H84").Value = 34908.668
$ws.Range("J84").Value = 34908.668
$ws.Range("L84").Value = 104726.004
$ws.Range("N84").Value = -114710.004

$ws.Range("H113").Value = 78081.94
$ws.Range("I113").Value = 40908.547
$ws.Range("K113").Value = 40908.547
$ws.Range("M113").Value = -38738.547

$ws.Range("H127").Value = 90000
$ws.Range("J127").Value = 90000
$ws.Range("L127").Value = 90000
$ws.Range("N127").Value = -99920

$ws.Range("H134").Value = 4709
$ws.Range("I134").Value = 4057.7058
$ws.Range("J134").Value = 8399.666999999999
$ws.Range("K134").Value = 12173.1174
$ws.Range("L134").Value = 25199.001
$ws.Range("M134").Value = -9638.117400000001
$ws.Range("N134").Value = -30269.001

$ws.Range("H136").Value = 5782
$ws.Range("I136").Value = 3912
$ws.Range("J136").Value = 6249.5
$ws.Range("K136").Value = 11736
$ws.Range("L136").Value = 18748.5
$ws.Range("M136").Value = -9186
$ws.Range("N136").Value = -23848.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 3858.3333
$ws.Range("I112").Value = 3830.2
$ws.Range("K112").Value = 11490.6
$ws.Range("M112").Value = -10382.6

$ws.Range("H126").Value = 17170.5
$ws.Range("I126").Value = 6665.3335
$ws.Range("K126").Value = 19996.0005
$ws.Range("M126").Value = -15056.0005

$ws.Range("H131").Value = 2219.44
$ws.Range("J131").Value = 2395.7
$ws.Range("L131").Value = 7187.099999999999
$ws.Range("N131").Value = -17267.1

$ws.Range("H132").Value = 4993.2
$ws.Range("I132").Value = 4993.6665
$ws.Range("J132").Value = 4992.5
$ws.Range("K132").Value = 44942.9985
$ws.Range("L132").Value = 44932.5
$ws.Range("M132").Value = -42412.9985
$ws.Range("N132").Value = -49992.5

$ws.Range("H140").Value = 3254.75
$ws.Range("I140").Value = 2704.389
$ws.Range("K140").Value = 8113.167
$ws.Range("M140").Value = -2933.167

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 25933.334
$ws.Range("J49").Value = 25933.334
$ws.Range("L49").Value = 25933.334
$ws.Range("N49").Value = -26301.334

$ws.Range("H112").Value = 89999
$ws.Range("J112").Value = 89999
$ws.Range("L112").Value = 89999
$ws.Range("N112").Value = -92215

$ws.Range("H122").Value = 3808.9092
$ws.Range("I122").Value = 2662.2222
$ws.Range("J122").Value = 8969
$ws.Range("K122").Value = 7986.6666
$ws.Range("L122").Value = 26907
$ws.Range("M122").Value = -5536.6666
$ws.Range("N122").Value = -31807

$ws.Range("H126").Value = 2803.7222
$ws.Range("I126").Value = 2104.1333
$ws.Range("K126").Value = 6312.3999
$ws.Range("M126").Value = -3842.3999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2191.5
$ws.Range("I22").Value = 1937.8125
$ws.Range("J22").Value = 2481.4285
$ws.Range("K22").Value = 1937.8125
$ws.Range("L22").Value = 2481.4285
$ws.Range("M22").Value = -1642.8125
$ws.Range("N22").Value = -3071.4285

$ws.Range("H27").Value = 2191.5
$ws.Range("I27").Value = 1937.8125
$ws.Range("J27").Value = 2481.4285
$ws.Range("K27").Value = 1937.8125
$ws.Range("L27").Value = 2481.4285
$ws.Range("M27").Value = -1830.8125
$ws.Range("N27").Value = -2695.4285

$ws.Range("H46").Value = 4069.4
$ws.Range("I46").Value = 4899
$ws.Range("J46").Value = 3516.3333
$ws.Range("K46").Value = 4899
$ws.Range("L46").Value = 3516.3333
$ws.Range("M46").Value = -4711
$ws.Range("N46").Value = -3892.3333

$ws.Range("H95").Value = 20479.6
$ws.Range("J95").Value = 20479.6
$ws.Range("L95").Value = 20479.6
$ws.Range("N95").Value = -25971.6

$ws.Range("H132").Value = 20898.795
$ws.Range("I132").Value = 23018.6
$ws.Range("K132").Value = 69055.79999999999
$ws.Range("M132").Value = -66525.79999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2966.35
$ws.Range("I122").Value = 2666.2856
$ws.Range("J122").Value = 3666.5
$ws.Range("K122").Value = 7998.8568
$ws.Range("L122").Value = 10999.5
$ws.Range("M122").Value = -5548.8568
$ws.Range("N122").Value = -15899.5

$ws.Range("H132").Value = 3073.0667
$ws.Range("I132").Value = 3216.923
$ws.Range("K132").Value = 9650.769
$ws.Range("M132").Value = -7120.769

$ws.Range("H136").Value = 33999
$ws.Range("J136").Value = 37497.5
$ws.Range("L136").Value = 112492.5
$ws.Range("N136").Value = -117592.5
